$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.499.55"
$ws.Range("E2").Value = "  -1.42%  "

$ws.Range("D3").Value = "2.344.80"
$ws.Range("E3").Value = "  -1.50%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.57"
$ws.Range("E5").Value = "  -2.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.25"
$ws.Range("E6").Value = "  +1.13%  "

$ws.Range("E7").Value = "  -1.77%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("E9").Value = "  -7.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.00"
$ws.Range("E10").Value = "  -1.20%  "

$ws.Range("E11").Value = "  -1.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.40"
$ws.Range("E12").Value = "  -2.18%  "

$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("E14").Value = "  -2.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.89"
$ws.Range("E15").Value = "  -5.91%  "

$ws.Range("D16").Value = "2.701.71"
$ws.Range("E16").Value = "  -1.20%  "

$ws.Range("D17").Value = "2.365.92"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").Value = "42.436.80"
$ws.Range("E18").Value = "  -1.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  -5.03%  "

$ws.Range("E20").Value = "  -2.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.85"
$ws.Range("E21").Value = "  -0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.57"
$ws.Range("E22").Value = "  +5.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.46"
$ws.Range("E23").Value = "  -7.51%  "

$ws.Range("E24").Value = "  -5.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.36"
$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.35"
$ws.Range("E27").Value = "  -2.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.76"
$ws.Range("E28").Value = "  -2.25%  "

$ws.Range("E29").Value = "  +2.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.97"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.23"
$ws.Range("E31").Value = "  -4.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0885"
$ws.Range("E32").Value = "  -3.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.05"
$ws.Range("E33").Value = "  +3.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.87"
$ws.Range("E34").Value = "  -8.49%  "

$ws.Range("E35").Value = "  +16.06%  "

$ws.Range("E37").Value = "  -5.56%  "

$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.85"
$ws.Range("E39").Value = "  -8.18%  "

$ws.Range("E40").Value = "  -6.10%  "

$ws.Range("E41").Value = "  +2.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.78"
$ws.Range("E42").Value = "  +1.73%  "

$ws.Range("E43").Value = "  -8.56%  "

$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.94"
$ws.Range("E45").Value = "  -4.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "112.23"
$ws.Range("E46").Value = "  -8.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.46"
$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.10"
$ws.Range("E48").Value = "  -10.48%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  -3.94%  "

$ws.Range("E50").Value = "  -2.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.86"
$ws.Range("E51").Value = "  +1.88%  "
